# Fruta / hortaliza, semanal
# Insert 3 new daily-price rows for Femacal de La Calera - Plátano (date 2021-11-11,
# serial 44511) ahead of the existing historical rows, pushing the rest of the
# table down by three rows (old row 432 -> new row 435, ... old row 526 -> new row 529).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 432:526 down to 435:529, opening up three blank rows at 432:434.
$ws.Rows("432:434").Insert()

# New row 432: Maduro
$ws.Cells.Item(432, 1).Value = 3
$ws.Cells.Item(432, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(432, 3).Value = "Coquimbo"
$ws.Cells.Item(432, 4).Value = 44511
$ws.Cells.Item(432, 5).Value = 5
$ws.Cells.Item(432, 6).Value = "Fruta"
$ws.Cells.Item(432, 7).Value = 100108
$ws.Cells.Item(432, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(432, 9).Value = 100108006
$ws.Cells.Item(432, 10).Value = "Plátano"
$ws.Cells.Item(432, 11).Value = "Sin especificar"
$ws.Cells.Item(432, 12).Value = "Maduro"
$ws.Cells.Item(432, 13).Value = 160
$ws.Cells.Item(432, 14).Value = 15000
$ws.Cells.Item(432, 15).Value = 15000
$ws.Cells.Item(432, 16).Value = 15000
$ws.Cells.Item(432, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(432, 18).Value = "Ecuador"
$ws.Cells.Item(432, 19).Value = 750
$ws.Cells.Item(432, 20).Value = 20

# New row 433: Pintón
$ws.Cells.Item(433, 1).Value = 3
$ws.Cells.Item(433, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(433, 3).Value = "Coquimbo"
$ws.Cells.Item(433, 4).Value = 44511
$ws.Cells.Item(433, 5).Value = 5
$ws.Cells.Item(433, 6).Value = "Fruta"
$ws.Cells.Item(433, 7).Value = 100108
$ws.Cells.Item(433, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(433, 9).Value = 100108006
$ws.Cells.Item(433, 10).Value = "Plátano"
$ws.Cells.Item(433, 11).Value = "Sin especificar"
$ws.Cells.Item(433, 12).Value = "Pintón"
$ws.Cells.Item(433, 13).Value = 360
$ws.Cells.Item(433, 14).Value = 17000
$ws.Cells.Item(433, 15).Value = 17000
$ws.Cells.Item(433, 16).Value = 17000
$ws.Cells.Item(433, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(433, 18).Value = "Ecuador"
$ws.Cells.Item(433, 19).Value = 850
$ws.Cells.Item(433, 20).Value = 20

# New row 434: Primera Pintón
$ws.Cells.Item(434, 1).Value = 3
$ws.Cells.Item(434, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(434, 3).Value = "Coquimbo"
$ws.Cells.Item(434, 4).Value = 44511
$ws.Cells.Item(434, 5).Value = 5
$ws.Cells.Item(434, 6).Value = "Fruta"
$ws.Cells.Item(434, 7).Value = 100108
$ws.Cells.Item(434, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(434, 9).Value = 100108006
$ws.Cells.Item(434, 10).Value = "Plátano"
$ws.Cells.Item(434, 11).Value = "Sin especificar"
$ws.Cells.Item(434, 12).Value = "Primera Pintón"
$ws.Cells.Item(434, 13).Value = 600
$ws.Cells.Item(434, 14).Value = 19000
$ws.Cells.Item(434, 15).Value = 20000
$ws.Cells.Item(434, 16).Value = 19533
$ws.Cells.Item(434, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(434, 18).Value = "Ecuador"
$ws.Cells.Item(434, 19).Value = 977
$ws.Cells.Item(434, 20).Value = 20
